$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("DATA")

# Insert a new worksheet "src" as the first sheet, containing a compact
# okres/potraty extract pulled by formula from the DATA sheet.
$src = $wb.Worksheets.Add()
$src.Move($wb.Worksheets.Item(1))
$src.Name = "src"

$src.Range("A1").Value = "okres"
$src.Range("B1").Value = "potraty"

$src.Range("A2").Formula = "=DATA!B9"
$src.Range("B2").Formula = "=DATA!C9"
$src.Range("A3").Formula = "=DATA!B11"
$src.Range("B3").Formula = "=DATA!C11"
$src.Range("A4").Formula = "=DATA!B12"
$src.Range("B4").Formula = "=DATA!C12"
$src.Range("A5").Formula = "=DATA!B13"
$src.Range("B5").Formula = "=DATA!C13"
$src.Range("A6").Formula = "=DATA!B14"
$src.Range("B6").Formula = "=DATA!C14"
$src.Range("A7").Formula = "=DATA!B15"
$src.Range("B7").Formula = "=DATA!C15"
$src.Range("A8").Formula = "=DATA!B16"
$src.Range("B8").Formula = "=DATA!C16"
$src.Range("A9").Formula = "=DATA!B17"
$src.Range("B9").Formula = "=DATA!C17"
$src.Range("A10").Formula = "=DATA!B18"
$src.Range("B10").Formula = "=DATA!C18"
$src.Range("A11").Formula = "=DATA!B19"
$src.Range("B11").Formula = "=DATA!C19"
$src.Range("A12").Formula = "=DATA!B20"
$src.Range("B12").Formula = "=DATA!C20"
$src.Range("A13").Formula = "=DATA!B21"
$src.Range("B13").Formula = "=DATA!C21"
$src.Range("A14").Formula = "=DATA!B22"
$src.Range("B14").Formula = "=DATA!C22"
$src.Range("A15").Formula = "=DATA!B24"
$src.Range("B15").Formula = "=DATA!C24"
$src.Range("A16").Formula = "=DATA!B25"
$src.Range("B16").Formula = "=DATA!C25"
$src.Range("A17").Formula = "=DATA!B26"
$src.Range("B17").Formula = "=DATA!C26"
$src.Range("A18").Formula = "=DATA!B27"
$src.Range("B18").Formula = "=DATA!C27"
$src.Range("A19").Formula = "=DATA!B28"
$src.Range("B19").Formula = "=DATA!C28"
$src.Range("A20").Formula = "=DATA!B29"
$src.Range("B20").Formula = "=DATA!C29"
$src.Range("A21").Formula = "=DATA!B30"
$src.Range("B21").Formula = "=DATA!C30"
$src.Range("A22").Formula = "=DATA!B32"
$src.Range("B22").Formula = "=DATA!C32"
$src.Range("A23").Formula = "=DATA!B33"
$src.Range("B23").Formula = "=DATA!C33"
$src.Range("A24").Formula = "=DATA!B34"
$src.Range("B24").Formula = "=DATA!C34"
$src.Range("A25").Formula = "=DATA!B35"
$src.Range("B25").Formula = "=DATA!C35"
$src.Range("A26").Formula = "=DATA!B36"
$src.Range("B26").Formula = "=DATA!C36"
$src.Range("A27").Formula = "=DATA!B37"
$src.Range("B27").Formula = "=DATA!C37"
$src.Range("A28").Formula = "=DATA!B38"
$src.Range("B28").Formula = "=DATA!C38"
$src.Range("A29").Formula = "=DATA!B40"
$src.Range("B29").Formula = "=DATA!C40"
$src.Range("A30").Formula = "=DATA!B41"
$src.Range("B30").Formula = "=DATA!C41"
$src.Range("A31").Formula = "=DATA!B42"
$src.Range("B31").Formula = "=DATA!C42"
$src.Range("A32").Formula = "=DATA!B44"
$src.Range("B32").Formula = "=DATA!C44"
$src.Range("A33").Formula = "=DATA!B45"
$src.Range("B33").Formula = "=DATA!C45"
$src.Range("A34").Formula = "=DATA!B46"
$src.Range("B34").Formula = "=DATA!C46"
$src.Range("A35").Formula = "=DATA!B47"
$src.Range("B35").Formula = "=DATA!C47"
$src.Range("A36").Formula = "=DATA!B48"
$src.Range("B36").Formula = "=DATA!C48"
$src.Range("A37").Formula = "=DATA!B49"
$src.Range("B37").Formula = "=DATA!C49"
$src.Range("A38").Formula = "=DATA!B50"
$src.Range("B38").Formula = "=DATA!C50"
$src.Range("A39").Formula = "=DATA!B52"
$src.Range("B39").Formula = "=DATA!C52"
$src.Range("A40").Formula = "=DATA!B53"
$src.Range("B40").Formula = "=DATA!C53"
$src.Range("A41").Formula = "=DATA!B54"
$src.Range("B41").Formula = "=DATA!C54"
$src.Range("A42").Formula = "=DATA!B55"
$src.Range("B42").Formula = "=DATA!C55"
$src.Range("A43").Formula = "=DATA!B57"
$src.Range("B43").Formula = "=DATA!C57"
$src.Range("A44").Formula = "=DATA!B58"
$src.Range("B44").Formula = "=DATA!C58"
$src.Range("A45").Formula = "=DATA!B59"
$src.Range("B45").Formula = "=DATA!C59"
$src.Range("A46").Formula = "=DATA!B60"
$src.Range("B46").Formula = "=DATA!C60"
$src.Range("A47").Formula = "=DATA!B61"
$src.Range("B47").Formula = "=DATA!C61"
$src.Range("A48").Formula = "=DATA!B63"
$src.Range("B48").Formula = "=DATA!C63"
$src.Range("A49").Formula = "=DATA!B64"
$src.Range("B49").Formula = "=DATA!C64"
$src.Range("A50").Formula = "=DATA!B65"
$src.Range("B50").Formula = "=DATA!C65"
$src.Range("A51").Formula = "=DATA!B66"
$src.Range("B51").Formula = "=DATA!C66"
$src.Range("A52").Formula = "=DATA!B68"
$src.Range("B52").Formula = "=DATA!C68"
$src.Range("A53").Formula = "=DATA!B69"
$src.Range("B53").Formula = "=DATA!C69"
$src.Range("A54").Formula = "=DATA!B70"
$src.Range("B54").Formula = "=DATA!C70"
$src.Range("A55").Formula = "=DATA!B71"
$src.Range("B55").Formula = "=DATA!C71"
$src.Range("A56").Formula = "=DATA!B72"
$src.Range("B56").Formula = "=DATA!C72"
$src.Range("A57").Formula = "=DATA!B74"
$src.Range("B57").Formula = "=DATA!C74"
$src.Range("A58").Formula = "=DATA!B75"
$src.Range("B58").Formula = "=DATA!C75"
$src.Range("A59").Formula = "=DATA!B76"
$src.Range("B59").Formula = "=DATA!C76"
$src.Range("A60").Formula = "=DATA!B77"
$src.Range("B60").Formula = "=DATA!C77"
$src.Range("A61").Formula = "=DATA!B78"
$src.Range("B61").Formula = "=DATA!C78"
$src.Range("A62").Formula = "=DATA!B79"
$src.Range("B62").Formula = "=DATA!C79"
$src.Range("A63").Formula = "=DATA!B80"
$src.Range("B63").Formula = "=DATA!C80"
$src.Range("A64").Formula = "=DATA!B82"
$src.Range("B64").Formula = "=DATA!C82"
$src.Range("A65").Formula = "=DATA!B83"
$src.Range("B65").Formula = "=DATA!C83"
$src.Range("A66").Formula = "=DATA!B84"
$src.Range("B66").Formula = "=DATA!C84"
$src.Range("A67").Formula = "=DATA!B85"
$src.Range("B67").Formula = "=DATA!C85"
$src.Range("A68").Formula = "=DATA!B86"
$src.Range("B68").Formula = "=DATA!C86"
$src.Range("A69").Formula = "=DATA!B88"
$src.Range("B69").Formula = "=DATA!C88"
$src.Range("A70").Formula = "=DATA!B89"
$src.Range("B70").Formula = "=DATA!C89"
$src.Range("A71").Formula = "=DATA!B90"
$src.Range("B71").Formula = "=DATA!C90"
$src.Range("A72").Formula = "=DATA!B91"
$src.Range("B72").Formula = "=DATA!C91"
$src.Range("A73").Formula = "=DATA!B93"
$src.Range("B73").Formula = "=DATA!C93"
$src.Range("A74").Formula = "=DATA!B94"
$src.Range("B74").Formula = "=DATA!C94"
$src.Range("A75").Formula = "=DATA!B95"
$src.Range("B75").Formula = "=DATA!C95"
$src.Range("A76").Formula = "=DATA!B96"
$src.Range("B76").Formula = "=DATA!C96"
$src.Range("A77").Formula = "=DATA!B97"
$src.Range("B77").Formula = "=DATA!C97"
$src.Range("A78").Formula = "=DATA!B98"
$src.Range("B78").Formula = "=DATA!C98"

# trailing blank rows (through row 85) to match the source layout
$src.Range("A79:B85").Value = $null

$src.Range("E4").Select()